# Actualización automática 2025-09-23 14:40:09
#
# A new advisor/client row ("RAMIREZ MOREIRA MAYRA JACQUELINE", under
# "OFICINA-CATAECSA") is inserted into the alphabetically-sorted client
# lists on the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, right before
# the existing "SALAZAR VERA ENRIQUE WILLIAM" row. All the rows that used
# to follow shift down by one, and the trailing "count" row (only present
# on "VENTAS POR GRUPO") has its "de 327" denominators bumped to "de 328"
# to reflect the extra data row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": columns A:R, new row lands at row 298.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$newRow1 = 298
$ws1.Rows.Item($newRow1).Insert()

$ws1.Cells.Item($newRow1, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item($newRow1, 2).Value = "RAMIREZ MOREIRA MAYRA JACQUELINE"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item($newRow1, $c).Value = 0
}

# The summary/count row used to be row 329 ("6 de 327", "25 de 327", ...);
# after the insert it is row 330 and the denominator needs to become 328.
$countRow1 = 330
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item($countRow1, $c)
    $cell.Value = $cell.Value2.Replace("de 327", "de 328")
}

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL": columns A:G, new row lands at row 302.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$newRow2 = 302
$ws2.Rows.Item($newRow2).Insert()

$ws2.Cells.Item($newRow2, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item($newRow2, 2).Value = "RAMIREZ MOREIRA MAYRA JACQUELINE"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item($newRow2, $c).Value = 0
}
